$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 134
$ws.Cells.Item(2, 4).Value = 58805.55
$ws.Cells.Item(2, 5).Value = 117.53

# Row 3
$ws.Cells.Item(3, 3).Value = 54
$ws.Cells.Item(3, 4).Value = 49155
$ws.Cells.Item(3, 5).Value = 970

# Row 4
$ws.Cells.Item(4, 3).Value = 55
$ws.Cells.Item(4, 4).Value = 46145

# Row 5
$ws.Cells.Item(5, 3).Value = 67
$ws.Cells.Item(5, 4).Value = 45573.91

# Row 6
$ws.Cells.Item(6, 3).Value = 62
$ws.Cells.Item(6, 4).Value = 44800
$ws.Cells.Item(6, 5).Value = 720

# Row 7
$ws.Cells.Item(7, 3).Value = 54
$ws.Cells.Item(7, 4).Value = 42870
$ws.Cells.Item(7, 5).Value = 785

# Row 8
$ws.Cells.Item(8, 3).Value = 67
$ws.Cells.Item(8, 4).Value = 39695
$ws.Cells.Item(8, 5).Value = 600

# Row 9
$ws.Cells.Item(9, 3).Value = 67
$ws.Cells.Item(9, 4).Value = 38585
$ws.Cells.Item(9, 5).Value = 555

# Row 10
$ws.Cells.Item(10, 3).Value = 67
$ws.Cells.Item(10, 4).Value = 33270
$ws.Cells.Item(10, 5).Value = 510

# Row 11
$ws.Cells.Item(11, 3).Value = 67
$ws.Cells.Item(11, 4).Value = 26840

# Row 12
$ws.Cells.Item(12, 3).Value = 67
$ws.Cells.Item(12, 4).Value = 24679.68
$ws.Cells.Item(12, 5).Value = 376.69

# Row 13
$ws.Cells.Item(13, 3).Value = 67
$ws.Cells.Item(13, 4).Value = 22148.31
$ws.Cells.Item(13, 5).Value = 335.74

# Row 14
$ws.Cells.Item(14, 3).Value = 67
$ws.Cells.Item(14, 4).Value = 14440.91
$ws.Cells.Item(14, 5).Value = 235.56

# Row 15
$ws.Cells.Item(15, 3).Value = 67
$ws.Cells.Item(15, 4).Value = 9869.65
$ws.Cells.Item(15, 5).Value = 156.38

# Row 16
$ws.Cells.Item(16, 3).Value = 67
$ws.Cells.Item(16, 4).Value = 8834.36
$ws.Cells.Item(16, 5).Value = 139.39

# Row 17
$ws.Cells.Item(17, 3).Value = 67
$ws.Cells.Item(17, 4).Value = 8112.31
$ws.Cells.Item(17, 5).Value = 121.76

# Row 18
$ws.Cells.Item(18, 3).Value = 67
$ws.Cells.Item(18, 4).Value = 7565.34
$ws.Cells.Item(18, 5).Value = 119.86

# Row 19
$ws.Cells.Item(19, 3).Value = 67
$ws.Cells.Item(19, 4).Value = 7490.67
$ws.Cells.Item(19, 5).Value = 117.35

# Row 20
$ws.Cells.Item(20, 3).Value = 67
$ws.Cells.Item(20, 4).Value = 7330.41
$ws.Cells.Item(20, 5).Value = 117.23

# Row 21
$ws.Cells.Item(21, 3).Value = 67
$ws.Cells.Item(21, 4).Value = 7222.86
$ws.Cells.Item(21, 5).Value = 113.43

# Row 22
$ws.Cells.Item(22, 3).Value = 67
$ws.Cells.Item(22, 4).Value = 7098.42
$ws.Cells.Item(22, 5).Value = 111.48

# Row 23
$ws.Cells.Item(23, 3).Value = 67
$ws.Cells.Item(23, 4).Value = 6706.38
$ws.Cells.Item(23, 5).Value = 96.19

# Row 24
$ws.Cells.Item(24, 3).Value = 67
$ws.Cells.Item(24, 4).Value = 6650.29
$ws.Cells.Item(24, 5).Value = 97.65

# Row 28
$ws.Cells.Item(28, 2).Value = 14
$ws.Cells.Item(28, 4).Value = 49.69
$ws.Cells.Item(28, 5).Value = 3.92

# Row 29
$ws.Cells.Item(29, 2).Value = 11
$ws.Cells.Item(29, 4).Value = 47.37
$ws.Cells.Item(29, 5).Value = 3.23

# Row 36
$ws.Cells.Item(36, 1).Value = "SAFCA CI (SAFC)"
$ws.Cells.Item(36, 2).Value = 8
$ws.Cells.Item(36, 3).Value = 4
$ws.Cells.Item(36, 4).Value = 22.3
$ws.Cells.Item(36, 5).Value = 1.27

# Row 37
$ws.Cells.Item(37, 1).Value = "PALM CI (PALC)"
$ws.Cells.Item(37, 3).Value = 5
$ws.Cells.Item(37, 4).Value = 21.55
$ws.Cells.Item(37, 5).Value = -4.17

# Row 38
$ws.Cells.Item(38, 1).Value = "CIE CI (CIEC)"
$ws.Cells.Item(38, 2).Value = 7
$ws.Cells.Item(38, 3).Value = 4
$ws.Cells.Item(38, 4).Value = 20.65
$ws.Cells.Item(38, 5).Value = -2.78

# Row 39
$ws.Cells.Item(39, 3).Value = 10
$ws.Cells.Item(39, 4).Value = 18.27
$ws.Cells.Item(39, 5).Value = -1.89

# Row 41
$ws.Cells.Item(41, 1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws.Cells.Item(41, 2).Value = 11
$ws.Cells.Item(41, 3).Value = 10
$ws.Cells.Item(41, 4).Value = 14.37
$ws.Cells.Item(41, 5).Value = 4.25

# Row 42
$ws.Cells.Item(42, 1).Value = "UNIWAX CI (UNXC)"
$ws.Cells.Item(42, 2).Value = 12
$ws.Cells.Item(42, 3).Value = 9
$ws.Cells.Item(42, 4).Value = 14.12
$ws.Cells.Item(42, 5).Value = 2.53

# Row 43
$ws.Cells.Item(43, 1).Value = "BANK OF AFRICA NG (BOAN)"
$ws.Cells.Item(43, 2).Value = 13
$ws.Cells.Item(43, 3).Value = 13
$ws.Cells.Item(43, 4).Value = 11.38
$ws.Cells.Item(43, 5).Value = -2.17

# Row 44
$ws.Cells.Item(44, 1).Value = "SUCRIVOIRE (SCRC)"
$ws.Cells.Item(44, 2).Value = 8
$ws.Cells.Item(44, 4).Value = 10.74
$ws.Cells.Item(44, 5).Value = -1.01

# Row 47
$ws.Cells.Item(47, 3).Value = 5
$ws.Cells.Item(47, 4).Value = 6.13
$ws.Cells.Item(47, 5).Value = -2.05

# Row 54
$ws.Cells.Item(54, 1).Value = "SOGB CI (SOGC)"
$ws.Cells.Item(54, 2).Value = 8
$ws.Cells.Item(54, 3).Value = 5
$ws.Cells.Item(54, 4).Value = 1.72
$ws.Cells.Item(54, 5).Value = 2.7

# Row 55
$ws.Cells.Item(55, 1).Value = "TOTAL"
$ws.Cells.Item(55, 2).Value = 0
$ws.Cells.Item(55, 3).Value = 66
$ws.Cells.Item(55, 4).Value = 0
$ws.Cells.Item(55, 5).Value = 0

# Row 56
$ws.Cells.Item(56, 1).Value = "SONATEL SN (SNTS)"
$ws.Cells.Item(56, 2).Value = 2
$ws.Cells.Item(56, 3).Value = 3
$ws.Cells.Item(56, 4).Value = -3.43
$ws.Cells.Item(56, 5).Value = 0.8

# Row 57
$ws.Cells.Item(57, 1).Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws.Cells.Item(57, 2).Value = 11
$ws.Cells.Item(57, 3).Value = 15
$ws.Cells.Item(57, 4).Value = -4.04
$ws.Cells.Item(57, 5).Value = -2.29

# Row 62
$ws.Cells.Item(62, 2).Value = 10
$ws.Cells.Item(62, 4).Value = -10.27
$ws.Cells.Item(62, 5).Value = 1.68
